$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two new columns before column D, shifting existing D:K to F:M
$ws.Range("D1:E1").EntireColumn.Insert()

# Copy number formats/styles from the (now-shifted) original D:E columns (now F:G)
# into the newly inserted D:E columns, per contiguous data block.
$ws.Range("F7:G35").Copy()
$ws.Range("D7").PasteSpecial(-4122)
$ws.Range("F38:G77").Copy()
$ws.Range("D38").PasteSpecial(-4122)
$ws.Range("F80:G102").Copy()
$ws.Range("D80").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new D:E columns with the latest two quarters of data
$newData = @{
  7 = @(43465, 43373)
  8 = @(2597700, 2986400)
  9 = @(2008500, 2383500)
  10 = @(589200, 602900)
  12 = @("NA", "NA")
  13 = @(0, 0)
  14 = @(210000, 0)
  15 = @(208800, 206300)
  17 = @(2678400, 2909700)
  18 = @(-80700, 76700)
  20 = @(-39700, -91800)
  21 = @(88400, 191300)
  22 = @(0, 0)
  23 = @(-120400, -15100)
  24 = @(-32200, -3900)
  25 = @(0, 0)
  26 = @(-88200, -11200)
  27 = @(-137200, -54000)
  28 = @(0, 0)
  29 = @(0, "NA")
  30 = @(0, 0)
  31 = @(0, 0)
  32 = @(39700, 91800)
  33 = @(-137200, -54000)
  34 = @(0, 0)
  35 = @(-137200, -54000)
  38 = @(43465, 43373)
  41 = @(232100, 203200)
  42 = @(0, 0)
  43 = @(865500, 1037200)
  44 = @(164700, 177900)
  45 = @(156600, 308100)
  46 = @(1418900, 1726400)
  47 = @(490500, 441500)
  48 = @(12928400, 12080900)
  49 = @(2029800, 2286200)
  50 = @(0, 0)
  51 = @(0, 0)
  52 = @(70600, 34100)
  53 = @(0, 0)
  54 = @(16938200, 16569100)
  57 = @(721900, 975500)
  58 = @(1027900, 290000)
  59 = @(1049000, 1290900)
  60 = @(2798800, 2556400)
  61 = @(5632400, 5689200)
  62 = @(790500, 815300)
  63 = @(0, 0)
  64 = @(0, 0)
  65 = @(0, 0)
  66 = @(10613100, 10249700)
  68 = @(0, 0)
  69 = @(0, 0)
  70 = @(245700, 238000)
  71 = @(0, 0)
  72 = @(-130400, -23800)
  73 = @(0, 0)
  74 = @(0, 0)
  75 = @(0, 0)
  76 = @(6079400, 6081400)
  77 = @(0, 0)
  80 = @(43465, 43373)
  81 = @(-137200, -54000)
  83 = @(208800, 206400)
  84 = @(0, 0)
  85 = @(0, 0)
  86 = @(0, 0)
  87 = @(0, 0)
  88 = @(0, 0)
  89 = @(229400, 370900)
  91 = @(-1081200, -871100)
  92 = @(0, 0)
  93 = @(0, 0)
  94 = @(-954100, -955900)
  96 = @(-231700, -231300)
  97 = @(0, 0)
  98 = @(0, 0)
  99 = @(0, 0)
  100 = @(753600, 506600)
  101 = @(0, 0)
  102 = @(28900, -78400)
}

foreach ($r in $newData.Keys) {
  $vals = $newData[$r]
  $ws.Cells.Item($r, 4).Value = $vals[0]
  $ws.Cells.Item($r, 5).Value = $vals[1]
}
